# Insert a new price-record row at row 106 (a new weekly "Ají" quote for
# Terminal Hortofrutícola Agro Chillán), pushing the existing rows 106-176
# down to 107-177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 106:176 down to 107:177, leaving a blank row 106.
$ws.Rows("106:106").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A106").Value = 7
$ws.Range("B106").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C106").Value = "Ñuble"
$ws.Range("D106").Value = 45016
$ws.Range("E106").Value = 16
$ws.Range("F106").Value = 100112021
$ws.Range("G106").Value = "Ají"
$ws.Range("H106").Value = "Cristal"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 50
$ws.Range("K106").Value = 15000
$ws.Range("L106").Value = 15000
$ws.Range("M106").Value = 15000
$ws.Range("N106").Value = "$/saco 25 kilos"
$ws.Range("O106").Value = "Región del Maule"
$ws.Range("P106").Value = 600
$ws.Range("Q106").Value = 25
$ws.Range("R106").Value = "Hortaliza"
